$wb = $excel.ActiveWorkbook

# --- Sheet "Cover": update the "Solution:" generated date ---
$wsCover = $wb.Worksheets.Item("Cover")
$wsCover.Range("C7").Value = "November 24, 2025"

# --- Sheet "Infrastructure Costs": rename "Software" category to "Software
#     Licenses", drop the Professional Services line items (rows 13-21),
#     and shrink the table / TOTAL row + formulas down to the new extent ---
$wsInfra = $wb.Worksheets.Item("Infrastructure Costs")
$wsInfra.Range("A7").Value = "Software Licenses"
$wsInfra.Range("A8").Value = "Software Licenses"
$wsInfra.Range("A9").Value = "Software Licenses"
$wsInfra.Range("A10").Value = "Software Licenses"
$wsInfra.Range("A11").Value = "Software Licenses"
$wsInfra.Rows("13:21").Delete()
$wsInfra.AutoFilterMode = $false
$wsInfra.Range("A2:K13").AutoFilter()

# --- Sheet "Credits": rename "Software" category to "Software Licenses"
#     and drop the Professional Services / Partner Credit row ---
$wsCredits = $wb.Worksheets.Item("Credits")
$wsCredits.Range("A4").Value = "Software Licenses"
$wsCredits.Rows(6).Delete()
$wsCredits.AutoFilterMode = $false
$wsCredits.Range("A2:D5").AutoFilter()

# --- Sheet "3-Year Summary": rename "Software" category to "Software
#     Licenses" and drop the Professional Services row (TOTAL moves up) ---
$wsSummary = $wb.Worksheets.Item("3-Year Summary")
$wsSummary.Range("A4").Value = "Software Licenses"
$wsSummary.Rows(6).Delete()
$wsSummary.AutoFilterMode = $false
$wsSummary.Range("A2:G6").AutoFilter()

# --- Workbook-level defined names: shrink the hidden _FilterDatabase
#     ranges to match the new table extents on the three affected sheets ---
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Infrastructure Costs!_FilterDatabase") {
    $n.RefersTo = "='Infrastructure Costs'!`$A`$2:`$K`$13"
  }
  elseif ($n.Name -eq "Credits!_FilterDatabase") {
    $n.RefersTo = "='Credits'!`$A`$2:`$D`$5"
  }
  elseif ($n.Name -eq "3-Year Summary!_FilterDatabase") {
    $n.RefersTo = "='3-Year Summary'!`$A`$2:`$G`$6"
  }
}
